$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated crypto price/volume data (and the Aptos/EnergySwap row swap)
# Each cell is written as literal text (NumberFormat "@") so values such as
# "1.001" or "29.894.49" are not reinterpreted as numbers/dates by Excel,
# then the style is reset to "Normal" so no stray formatting is introduced.
function Set-TextCell($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextCell "D2" "29.894.49"
Set-TextCell "E2" "  +0.54%  "

Set-TextCell "D3" "1.895.75"
Set-TextCell "E3" "  +0.53%  "

Set-TextCell "E4" "  +0.10%  "

Set-TextCell "D5" "0.7830"
Set-TextCell "E5" "  -1.32%  "

Set-TextCell "D6" "243.75"
Set-TextCell "E6" "  +1.07%  "

Set-TextCell "D7" "1.001"
Set-TextCell "E7" "  +0.13%  "

Set-TextCell "E8" "  -0.84%  "

Set-TextCell "D9" "25.72"
Set-TextCell "E9" "  +1.08%  "

Set-TextCell "D10" "0.07265"
Set-TextCell "E10" "  +3.84%  "

Set-TextCell "D11" "0.08110"
Set-TextCell "E11" "  +0.92%  "

Set-TextCell "D12" "0.7742"
Set-TextCell "E12" "  +1.85%  "

Set-TextCell "D13" "5.477"
Set-TextCell "E13" "  +3.54%  "

Set-TextCell "D14" "1.896.49"
Set-TextCell "E14" "  +1.36%  "

Set-TextCell "D15" "94.49"
Set-TextCell "E15" "  +2.56%  "

Set-TextCell "D16" "6.215"
Set-TextCell "E16" "  +4.88%  "

Set-TextCell "D17" "29.896.39"

Set-TextCell "E18" "  +0.96%  "

Set-TextCell "D19" "246.11"
Set-TextCell "E19" "  +1.22%  "

Set-TextCell "D20" "0.000007840"
Set-TextCell "E20" "  +2.24%  "

Set-TextCell "D21" "1.001"
Set-TextCell "E21" "  +0.10%  "

Set-TextCell "D22" "8.143"
Set-TextCell "E22" "  -0.25%  "

Set-TextCell "D23" "2.131.85"
Set-TextCell "E23" "  +0.98%  "

Set-TextCell "E24" "  +0.06%  "

Set-TextCell "D25" "0.1599"
Set-TextCell "E25" "  -4.39%  "

Set-TextCell "D26" "9.465"
Set-TextCell "E26" "  +2.03%  "

Set-TextCell "D27" "164.38"
Set-TextCell "E27" "  +0.60%  "

Set-TextCell "D28" "18.79"
Set-TextCell "E28" "  +0.98%  "

Set-TextCell "D29" "2.022"
Set-TextCell "E29" "  -1.22%  "

Set-TextCell "D30" "1.435"
Set-TextCell "E30" "  +3.48%  "

Set-TextCell "D31" "1.544"
Set-TextCell "E31" "  +0.91%  "

Set-TextCell "D32" "4.481"
Set-TextCell "E32" "  +2.60%  "

Set-TextCell "D33" "0.05583"
Set-TextCell "E33" "  -1.61%  "

Set-TextCell "D34" "4.081"
Set-TextCell "E34" "  +0.87%  "

Set-TextCell "E35" "  -1.11%  "

Set-TextCell "E36" "  +2.95%  "

Set-TextCell "D37" "1.003"
Set-TextCell "E37" "  +0.66%  "

Set-TextCell "E38" "  +2.49%  "

Set-TextCell "D39" "0.01935"
Set-TextCell "E39" "  +1.89%  "

Set-TextCell "D40" "2.786"
Set-TextCell "E40" "  +0.62%  "

Set-TextCell "D41" "1.145.38"
Set-TextCell "E41" "  +12.04%  "

Set-TextCell "D42" "0.4460"
Set-TextCell "E42" "  +1.46%  "

Set-TextCell "D43" "74.04"
Set-TextCell "E43" "  +2.56%  "

Set-TextCell "D44" "5.954"
Set-TextCell "E44" "  +2.40%  "

Set-TextCell "D45" "0.8529"
Set-TextCell "E45" "  +2.10%  "

Set-TextCell "D46" "1.001"
Set-TextCell "E46" "  +0.10%  "

Set-TextCell "D47" "1.892"
Set-TextCell "E47" "  +1.62%  "

Set-TextCell "D48" "3.137"
Set-TextCell "E48" "  +8.12%  "

Set-TextCell "D49" "102.15"
Set-TextCell "E49" "  -0.44%  "

Set-TextCell "B50" "Aptos"
Set-TextCell "C50" "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextCell "D50" "7.534"
Set-TextCell "E50" "  +1.65%  "

Set-TextCell "B51" "EnergySwap"
Set-TextCell "C51" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextCell "D51" "9.738"
Set-TextCell "E51" "  -1.13%  "
